# DCF Update for Cruise Liners
# Added other statistical non-GAAP measures: Ticket APCD, Extras PCD, Extras APCD
#
# Model sheet, rows 74-82 (the "RCL ..." benchmarking block) are restructured:
#   - new row 74  "Extras per PCD"   = row4 / row66   (was blank)
#   - row 75      cleared (was "RCL Cash extraction per PCD" raw data)
#   - row 76      "Ticket per APCD"  = row3 / row$67  (was "RCL Cash extraction per APCD" raw data)
#   - row 77      "Extras per APCD" = row4 / row$67  (was a blank spacer row)
#   - row 78      cleared (was "RCL/NCLH Cash extract per PCD" ratio formulas)
#   - row 79      cleared (was "RCL/NCLH Cash extract per APCD" ratio formulas)
#   - row 81      cleared (was "RCL Ticket per PCD" raw data)
#   - row 82      cleared (was "RCL/NCLH ticket price per APCD" ratio formulas)

$wb = $excel.ActiveWorkbook
$wsModel = $wb.Worksheets.Item("Model")
$wsMain = $wb.Worksheets.Item("Main")

# Columns that carry data in this block: E:I and O:AA
$dataCols = @("E","F","G","H","I","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# --- Wipe out the rows that are going away / being replaced ---------------
$wsModel.Range("B75:AA75").ClearContents()
$wsModel.Range("B78:AA78").ClearContents()
$wsModel.Range("B79:AA79").ClearContents()
$wsModel.Range("B81:AA81").ClearContents()
$wsModel.Range("B82:AA82").ClearContents()

# --- Row 74: Extras per PCD = (Onboard extras row4) / (Passenger cruise days row66)
$wsModel.Range("B74").Value = "Extras per PCD"
foreach ($col in $dataCols) {
    $wsModel.Range($col + "74").Formula = "=" + $col + "4/" + $col + "66"
}

# --- Row 76: Ticket per APCD = (Passenger tickets row3) / (Capacity days row$67)
$wsModel.Range("B76").Value = "Ticket per APCD"
foreach ($col in $dataCols) {
    $wsModel.Range($col + "76").Formula = "=" + $col + "3/" + $col + "`$67"
}

# --- Row 77: Extras per APCD = (Onboard extras row4) / (Capacity days row$67)
$wsModel.Range("B77").Value = "Extras per APCD"
foreach ($col in $dataCols) {
    $wsModel.Range($col + "77").Formula = "=" + $col + "4/" + $col + "`$67"
}

# --- View state -------------------------------------------------------------
# Model: pane scrolled so the bottom-right pane shows S5, no longer the active tab
$wsModel.Activate()
$wsModel.Range("S5").Select()

# Main: now the active tab, with I18 selected
$wsMain.Activate()
$wsMain.Range("I18").Select()
